$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move together with each record (row).
# Columns A,B,C,E,F,G,H,I,J,K are identical across all rows and stay untouched.
$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)   # D, L, M, N, O, P, Q, R, S, T

# Snapshot current values (rows 2..18) for each of those columns.
$snapshot = @{}
for ($r = 2; $r -le 18; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Row r (destination) receives the former contents of row $mapping[r] (source).
$mapping = @{
    2  = 5
    3  = 8
    4  = 10
    5  = 13
    6  = 2
    7  = 17
    8  = 12
    9  = 14
    10 = 3
    11 = 4
    12 = 9
    13 = 6
    14 = 7
    15 = 11
    16 = 18
    17 = 15
    18 = 16
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
    }
}
